# Apply cryptos list update (price/volume refresh + two row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.080.50"
$ws.Cells.Item(2, 5).Value = "  -0.97%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.651.89"
$ws.Cells.Item(3, 5).Value = "  -1.04%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.005"
$ws.Cells.Item(4, 5).Value = "  -0.58%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "218.69"
$ws.Cells.Item(5, 5).Value = "  -0.77%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.5255"
$ws.Cells.Item(6, 5).Value = "  -1.00%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.005"
$ws.Cells.Item(7, 5).Value = "  -0.54%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2672"
$ws.Cells.Item(8, 5).Value = "  +0.71%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06365"
$ws.Cells.Item(9, 5).Value = "  -0.14%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "20.51"
$ws.Cells.Item(10, 5).Value = "  -2.11%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07695"
$ws.Cells.Item(11, 5).Value = "  -2.12%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "4.593"
$ws.Cells.Item(12, 5).Value = "  +1.35%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.663.91"
$ws.Cells.Item(13, 5).Value = "  -0.45%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "1.880.72"
$ws.Cells.Item(14, 5).Value = "  -0.98%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.5612"
$ws.Cells.Item(15, 5).Value = "  -0.11%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.0₅8224"
$ws.Cells.Item(16, 5).Value = "  +1.05%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "65.58"
$ws.Cells.Item(17, 5).Value = "  -0.52%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "26.094.89"
$ws.Cells.Item(18, 5).Value = "  -1.07%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.57%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.679"
$ws.Cells.Item(20, 5).Value = "  -0.97%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "190.37"
$ws.Cells.Item(21, 5).Value = "  -5.20%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "10.30"
$ws.Cells.Item(22, 5).Value = "  -0.04%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.963"
$ws.Cells.Item(23, 5).Value = "  -1.63%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.006"
$ws.Cells.Item(24, 5).Value = "  -0.57%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.74%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.1200"
$ws.Cells.Item(26, 5).Value = "  -1.24%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.244"
$ws.Cells.Item(27, 5).Value = "  -0.27%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -1.71%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.520"
$ws.Cells.Item(29, 5).Value = "  +0.46%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.05640"
$ws.Cells.Item(30, 5).Value = "  -4.39%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.270"
$ws.Cells.Item(31, 5).Value = "  -1.12%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.487"
$ws.Cells.Item(32, 5).Value = "  -0.96%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.374"
$ws.Cells.Item(33, 5).Value = "  +1.33%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.578"
$ws.Cells.Item(34, 5).Value = "  -1.50%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.801"
$ws.Cells.Item(35, 5).Value = "  -1.06%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "HuobiToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.411"
$ws.Cells.Item(36, 5).Value = "  -0.90%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "ARBITRUM"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.9454"
$ws.Cells.Item(37, 5).Value = "  -2.15%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.5768"
$ws.Cells.Item(38, 5).Value = "  -0.60%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01590"
$ws.Cells.Item(39, 5).Value = "  -1.81%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "5.963"
$ws.Cells.Item(40, 5).Value = "  -0.11%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -1.80%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -0.58%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.021.37"
$ws.Cells.Item(43, 5).Value = "  -5.19%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "101.19"
$ws.Cells.Item(44, 5).Value = "  -1.88%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.791.69"
$ws.Cells.Item(45, 5).Value = "  -1.01%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "58.29"
$ws.Cells.Item(46, 5).Value = "  -0.42%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.002"
$ws.Cells.Item(47, 5).Value = "  -1.28%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +3.49%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.4346"
$ws.Cells.Item(49, 5).Value = "  -1.68%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "8.001"
$ws.Cells.Item(50, 5).Value = "  -0.30%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -3.83%  "
